# CryCompanywiseStockReport_1.xlsx update
# Applies the per-row quantity (F) / stock value (G) recalculations,
# two pairs of duplicate-product rows that get their data swapped
# (HIM Total Care Baby Pants Diapers, HUL Kissan jam / Knorr / pears,
# KUS Floor Wiper, Rasna Nagpur Orange), and the resulting Sub Total /
# Grand Total rollups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F32").Value = 32
$ws.Range("G32").Value = 1547.52
$ws.Range("B34").Value = 61899.66
$ws.Range("F36").Value = 99
$ws.Range("G36").Value = 19480.23
$ws.Range("F41").Value = 226
$ws.Range("G41").Value = 43593.14
$ws.Range("F42").Value = 70
$ws.Range("G42").Value = 2948.4
$ws.Range("F45").Value = 86
$ws.Range("G45").Value = 1955.64
$ws.Range("F48").Value = 243
$ws.Range("G48").Value = 13632.3
$ws.Range("F51").Value = 157
$ws.Range("G51").Value = 14685.78
$ws.Range("F52").Value = 65
$ws.Range("G52").Value = 3835
$ws.Range("F53").Value = 38
$ws.Range("G53").Value = 623.58
$ws.Range("B66").Value = 220809.95
$ws.Range("F140").Value = 99
$ws.Range("G140").Value = 5206.41
$ws.Range("B147").Value = 23027.62
$ws.Range("F175").Value = 31
$ws.Range("G175").Value = 8990.309999999999
$ws.Range("F182").Value = 28
$ws.Range("G182").Value = 2507.12
$ws.Range("F186").Value = 31
$ws.Range("G186").Value = 1341.68
$ws.Range("B193").Value = 69627.53
$ws.Range("F222").Value = 1086
$ws.Range("G222").Value = 20091
$ws.Range("B229").Value = 33960.37
$ws.Range("F268").Value = 21
$ws.Range("G268").Value = 2671.41
$ws.Range("F277").Value = 14
$ws.Range("G277").Value = 705.74
$ws.Range("F278").Value = 50
$ws.Range("G278").Value = 6777
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 27
$ws.Range("G291").Value = 2313.36
$ws.Range("B295").Value = 132314.7
$ws.Range("B304").Value = 63520
$ws.Range("E304").Value = 153.4
$ws.Range("F304").Value = 39
$ws.Range("G304").Value = 5626.92
$ws.Range("B305").Value = 55373
$ws.Range("E305").Value = 163.62
$ws.Range("F305").Value = -94
$ws.Range("G305").Value = -13562.32
$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 30
$ws.Range("G306").Value = 4304.4
$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 80
$ws.Range("G308").Value = 3811.2
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12
$ws.Range("F314").Value = 1
$ws.Range("G314").Value = 521.64
$ws.Range("B317").Value = 60325
$ws.Range("E317").Value = 151.57
$ws.Range("F317").Value = -102
$ws.Range("G317").Value = -12939.72
$ws.Range("B318").Value = 63560
$ws.Range("E318").Value = 134.87
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 126.86
$ws.Range("B328").Value = 2488.03
$ws.Range("F352").Value = 128
$ws.Range("G352").Value = 15988.48
$ws.Range("B356").Value = 80302.69
$ws.Range("F361").Value = 259
$ws.Range("G361").Value = 36412.81
$ws.Range("B363").Value = 81666.92999999999
$ws.Range("F368").Value = 63
$ws.Range("G368").Value = 2017.26
$ws.Range("B372").Value = 66353.67
$ws.Range("B381").Value = 47097
$ws.Range("D381").Value = 112.28
$ws.Range("E381").Value = 134.16
$ws.Range("F381").Value = 15
$ws.Range("G381").Value = 1684.2
$ws.Range("B382").Value = 58047
$ws.Range("D382").Value = 105.54
$ws.Range("E382").Value = 126.1
$ws.Range("F382").Value = 32
$ws.Range("G382").Value = 3377.28
$ws.Range("F387").Value = 454
$ws.Range("G387").Value = 43856.4
$ws.Range("B389").Value = 61062.69
$ws.Range("F392").Value = 148
$ws.Range("G392").Value = 5538.16
$ws.Range("F396").Value = 140
$ws.Range("G396").Value = 3567.2
$ws.Range("F399").Value = 274
$ws.Range("G399").Value = 26920.5
$ws.Range("F402").Value = 61
$ws.Range("G402").Value = 2092.91
$ws.Range("F406").Value = 79
$ws.Range("G406").Value = 8126.73
$ws.Range("F415").Value = 62
$ws.Range("G415").Value = 3391.4
$ws.Range("F416").Value = 77
$ws.Range("G416").Value = 2260.72
$ws.Range("B417").Value = 179709.61
$ws.Range("F430").Value = 237
$ws.Range("G430").Value = 10968.36
$ws.Range("B438").Value = 27227.27
$ws.Range("F454").Value = 82
$ws.Range("G454").Value = 23202.72
$ws.Range("F455").Value = 47
$ws.Range("G455").Value = 10444.81
$ws.Range("B458").Value = 103658.46
$ws.Range("F467").Value = 17
$ws.Range("G467").Value = 11247.2
$ws.Range("B476").Value = 51042.49
$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 86
$ws.Range("G506").Value = 2823.38
$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79
$ws.Range("B508").Value = 41849.99
$ws.Range("F520").Value = 35
$ws.Range("G520").Value = 959
$ws.Range("B525").Value = 132479.73
$ws.Range("F558").Value = 222
$ws.Range("G558").Value = 27050.7
$ws.Range("B561").Value = 31862.42
$ws.Range("F620").Value = 378
$ws.Range("G620").Value = 29707.02
$ws.Range("F625").Value = 340
$ws.Range("G625").Value = 12522.2
$ws.Range("B628").Value = 225564.19
$ws.Range("F659").Value = 45
$ws.Range("G659").Value = 2409.3
$ws.Range("F662").Value = 51
$ws.Range("G662").Value = 4095.81
$ws.Range("B668").Value = 13744.81
$ws.Range("F674").Value = 957
$ws.Range("G674").Value = 156096.27
$ws.Range("B680").Value = 157108.82
$ws.Range("F708").Value = 131
$ws.Range("G708").Value = 4831.28
$ws.Range("B713").Value = 72237.2
$ws.Range("B718").Value = 2982892.15
$ws.Range("B719").Value = 2982892.15
